{"js": "// Fill in the first empty data row of the time-tracking table with a\n// new entry: date \"16.3\", hours \"3\", and description\n// \"UML kuvaaja ja ohjelman hienos\u00e4\u00e4t\u00f6\u00e4\".\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.rows.load(\"items\");\nawait context.sync();\n\n// Locate the first row (after the header) whose first cell has no text \u2014\n// i.e. the next free entry row in the log.\nlet targetRowIndex = -1;\nconst rows = table.rows.items;\nfor (let i = 1; i < rows.length; i++) {\n  const cell = table.getCell(i, 0);\n  cell.body.load(\"text\");\n  await context.sync();\n  if (cell.body.text.trim() === \"\") {\n    targetRowIndex = i;\n    break;\n  }\n}\n\nif (targetRowIndex === -1) {\n  throw new Error(\"No empty row found to fill in.\");\n}\n\nconst newValues = [\"16.3\", \"3\", \"UML kuvaaja ja ohjelman hienos\u00e4\u00e4t\u00f6\u00e4\"];\n\nfor (let col = 0; col < newValues.length; col++) {\n  const cell = table.getCell(targetRowIndex, col);\n  cell.body.paragraphs.load(\"items\");\n  await context.sync();\n  const paragraph = cell.body.paragraphs.items[0];\n  paragraph.insertText(newValues[col], Word.InsertLocation.start);\n}\n\nawait context.sync();\n", "ps1": "# Fill in the first empty data row of the time-tracking table with a\n# new entry: date \"16.3\", hours \"3\", and description\n# \"UML kuvaaja ja ohjelman hienos\u00e4\u00e4t\u00f6\u00e4\".\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Find the first row after the header whose first cell is still empty \u2014\n# i.e. the next free entry row in the log.\n$targetRow = 0\nfor ($i = 2; $i -le $t.Rows.Count; $i++) {\n    $cellText = $t.Rows.Item($i).Cells.Item(1).Range.Text -replace '[\\r\\a]', ''\n    if ($cellText -eq \"\") {\n        $targetRow = $i\n        break\n    }\n}\n\nif ($targetRow -eq 0) {\n    throw \"No empty row found to fill in.\"\n}\n\n$row = $t.Rows.Item($targetRow)\n$row.Cells.Item(1).Range.Text = \"16.3\"\n$row.Cells.Item(2).Range.Text = \"3\"\n$row.Cells.Item(3).Range.Text = \"UML kuvaaja ja ohjelman hienos\u00e4\u00e4t\u00f6\u00e4\"\n"}
